# Update "想去人数" (column F) figures across the 展览, 演出 and 全部类型
# sheets to reflect the latest generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$expoChanges = @{
    2  = 3105
    3  = 508
    4  = 71
    5  = 68
    8  = 22
    9  = 1092
    10 = 15230
    11 = 210
    12 = 154
    13 = 700
    14 = 6034
    15 = 616
    16 = 94
    17 = 60
    18 = 97
    19 = 1255
    20 = 25
    21 = 108
    22 = 7
    23 = 206
    24 = 841
    25 = 8
    27 = 126
    28 = 10889
    29 = 1222
    31 = 100
    32 = 144
    34 = 255
}
foreach ($row in $expoChanges.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoChanges[$row]
}

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(3, 6).Value = 17

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$allChanges = @{
    3  = 3105
    4  = 508
    5  = 71
    6  = 68
    9  = 22
    10 = 1092
    11 = 15230
    12 = 210
    13 = 154
    14 = 701
    15 = 6034
    16 = 616
    17 = 94
    18 = 60
    19 = 97
    20 = 1255
    21 = 25
    22 = 108
    23 = 7
    24 = 206
    25 = 841
    26 = 8
    28 = 126
    29 = 17
    30 = 10889
    31 = 1222
    33 = 100
    34 = 144
    36 = 255
}
foreach ($row in $allChanges.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allChanges[$row]
}
